$d = $word.ActiveDocument

# Locate the start of the section to rewrite: "collects data on adult salmonids on Battle Creek"
$findRange = $d.Content
$findRange.Find.Execute("collects data on adult salmonids on Battle Creek", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$startPos = $findRange.Start

# The paragraph containing this text; its Range.End includes the trailing paragraph mark,
# so back up one character to land on the last real character of the paragraph.
$para = $findRange.Paragraphs(1)
$paraContentEnd = $para.Range.End - 1

# Range spanning from the start of "collects data..." through the end of the paragraph
# (i.e. through "...Sacramento River Watershed."). We rebuild this whole span as a fresh
# run-for-run OOXML fragment so the new wording lands in the exact run layout the edit calls for,
# while still reproducing (unchanged) the trailing runs that follow the rewritten sentences.
$targetRange = $d.Range($startPos, $paraContentEnd)

$xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/></w:rPr><w:t xml:space="preserve">collects data on adult salmonids on Battle Creek. </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/></w:rPr><w:t xml:space="preserve">Data are collected via snorkel surveys, video camera systems, and a trap at the spawning building. Snorkel surveys are conducted annually, and video camera systems operate 24 hours a day, 7 days a week when the trap is not operating due to temperature constraints. </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/></w:rPr><w:t xml:space="preserve">Data from this monitoring </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/></w:rPr><w:t>are</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/></w:rPr><w:t xml:space="preserve"> used to estimate adult escapement (upstream passage) abundance and timing, spawner abundance, and other important metrics for adult salmonids in the watershed. These</w:t></w:r><w:r w:rsidRPr="0019519D"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/></w:rPr><w:t xml:space="preserve"> data will also be used to inform the development of a juvenile production estimate (JPE) for spring-run Chinook salmon in the Sacramento River Watershed</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/></w:rPr><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$targetRange.InsertXML($xml)
